{"js": "// Update the title's academic-year suffix from \" (2023/24)\" to \" (2024/25)\".\n// The target text lives in its own run (with rFonts Avenir Next, bold) at the\n// end of the document's first (title) paragraph, e.g.:\n//   \"Data Analysis Project \u2013 GIS Mapping Guidance Sheet (2023/24)\"\n// becomes:\n//   \"Data Analysis Project \u2013 GIS Mapping Guidance Sheet (2024/25)\"\nconst oldText = \" (2023/24)\";\nconst newText = \" (2024/25)\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // \"Replace\" keeps the matched range's existing character formatting\n  // (font / bold / etc.) while swapping in the new text.\n  results.items[i].insertText(newText, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the title's academic-year suffix from \" (2023/24)\" to \" (2024/25)\".\n# The target text lives in its own run (rFonts Avenir Next, bold) at the end\n# of the document's first (title) paragraph, e.g.:\n#   \"Data Analysis Project - GIS Mapping Guidance Sheet (2023/24)\"\n# becomes:\n#   \"Data Analysis Project - GIS Mapping Guidance Sheet (2024/25)\"\n\n$d = $word.ActiveDocument\n$oldText = \" (2023/24)\"\n$newText = \" (2024/25)\"\n\n# Bounded loop: there is a single occurrence in this document, but the loop\n# stays safe even if the text appeared more than once (no infinite loop risk).\nfor ($i = 0; $i -lt 50; $i++) {\n  $rng = $d.Content\n  $find = $rng.Find\n  $find.Text = $oldText\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.Forward = $true\n  if (-not $find.Execute()) { break }\n\n  # Capture the matched run's formatting so the replacement keeps it,\n  # then delete the old text and insert the new text in its place.\n  $fontName = $rng.Font.Name\n  $isBold = $rng.Font.Bold\n  $start = $rng.Start\n  $rng.Delete()\n\n  $ip = $d.Range($start, $start)\n  $ip.InsertAfter($newText)\n\n  $newRng = $d.Range($start, $start + $newText.Length)\n  $newRng.Font.Name = $fontName\n  $newRng.Font.Bold = $isBold\n  $newRng.Font.BoldBi = $isBold\n}\n"}
